# Auto-generated Excel COM-interop script applying the cryptos.xlsx price/volume update
# described by the commit "Updated cryptos list on Sat Mar 18 14:44:22 UTC 2023 with GitHub Actions".
#
# The sheet stores every Price/Volume cell as literal text (t="inlineStr" in the OOXML,
# no numeric formatting applied). Assigning a clean-looking numeric string straight to
# .Value would let Excel auto-convert it to a real number, so any Price cell whose text
# parses as a plain number is entered with a leading apostrophe (forces text entry, the
# same trick a human would use in the Excel UI) and then its style is reset back to
# "Normal" so no stray quote-prefix/number-format style is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.597.06'
$ws.Range("E2").Value = '  +4.58%  '

$ws.Range("D3").Value = '1.824.09'
$ws.Range("E3").Value = '  +5.73%  '

$ws.Range("E4").Value = '  -0.70%  '

$ws.Range("D5").Value = "'338.12"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.61%  '

$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.36%  '

$ws.Range("D7").Value = "'0.3827"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +2.40%  '

$ws.Range("D8").Value = "'0.3537"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +5.54%  '

$ws.Range("D9").Value = "'49.93"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.37%  '

$ws.Range("D10").Value = "'1.239"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.74%  '

$ws.Range("D11").Value = "'0.07744"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.34%  '

$ws.Range("D12").Value = "'1.002"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.79%  '

$ws.Range("D13").Value = "'22.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +11.14%  '

$ws.Range("D14").Value = "'6.634"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.32%  '

$ws.Range("D15").Value = '1.826.80'
$ws.Range("E15").Value = '  +5.35%  '

$ws.Range("D16").Value = "'7.191"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.66%  '

$ws.Range("D17").Value = "'0.00001125"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +5.20%  '

$ws.Range("D18").Value = "'0.06723"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.39%  '

$ws.Range("D19").Value = "'87.24"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +5.68%  '

$ws.Range("D20").Value = "'1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.33%  '

$ws.Range("D21").Value = "'17.67"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +7.05%  '

$ws.Range("D22").Value = "'6.551"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +7.16%  '

$ws.Range("D23").Value = "'13.15"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.53%  '

$ws.Range("D24").Value = '27.608.69'
$ws.Range("E24").Value = '  +4.39%  '

$ws.Range("E25").Value = '  +0.82%  '

$ws.Range("D26").Value = "'2.667"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +11.99%  '

$ws.Range("D27").Value = "'22.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +14.23%  '

$ws.Range("D28").Value = "'1.485"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.67%  '

$ws.Range("D29").Value = "'152.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.47%  '

$ws.Range("D30").Value = '2.031.30'
$ws.Range("E30").Value = '  +5.37%  '

$ws.Range("D31").Value = "'135.14"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.25%  '

$ws.Range("D32").Value = "'6.345"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.55%  '

$ws.Range("D33").Value = "'4.086"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.46%  '

$ws.Range("D34").Value = "'13.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +9.74%  '

$ws.Range("D35").Value = "'0.08812"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.64%  '

$ws.Range("D36").Value = "'1.699"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.01%  '

$ws.Range("D37").Value = "'5.643"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +5.09%  '

$ws.Range("D38").Value = "'0.7035"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +14.35%  '

$ws.Range("D39").Value = "'9.150"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +7.85%  '

$ws.Range("B40").Value = 'Hedera'
$ws.Range("C40").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D40").Value = "'0.06532"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.18%  '

$ws.Range("B41").Value = 'Algorand'
$ws.Range("C41").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D41").Value = "'0.2262"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.92%  '

$ws.Range("D42").Value = "'0.02405"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.33%  '

$ws.Range("D43").Value = "'1.300"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.14%  '

$ws.Range("D44").Value = "'14.82"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +6.25%  '

$ws.Range("D45").Value = "'0.6613"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +11.10%  '

$ws.Range("D46").Value = "'1.0000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.48%  '

$ws.Range("D47").Value = "'3.905"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.02%  '

$ws.Range("D48").Value = "'2.187"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.60%  '

$ws.Range("D49").Value = "'133.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.76%  '

$ws.Range("D50").Value = "'0.07308"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.62%  '

$ws.Range("D51").Value = "'81.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.21%  '

